# Updated cryptos list on Mon Aug 14 04:32:10 UTC 2023 with GitHub Actions
#
# Writes string values into cells while forcing "Text" interpretation so that
# numeric-looking strings (prices like "0.9986") are preserved verbatim as
# text instead of being auto-converted into numbers by Excel, matching the
# original workbook's inline-string cell layout.

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "29.391.10"
Set-TextValue $ws.Range("E2") "  -0.07%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.846.15"
Set-TextValue $ws.Range("E3") "  -0.22%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "0.9986"
Set-TextValue $ws.Range("E4") "  -0.17%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "240.77"
Set-TextValue $ws.Range("E5") "  +0.07%  "

# Row 6 - XRP
Set-TextValue $ws.Range("E6") "  +0.22%  "

# Row 7 - USDC
Set-TextValue $ws.Range("E7") "  -0.10%  "

# Row 8 - Dogecoin
Set-TextValue $ws.Range("D8") "0.07489"
Set-TextValue $ws.Range("E8") "  -2.08%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.2908"
Set-TextValue $ws.Range("E9") "  +0.03%  "

# Row 10 - Solana
Set-TextValue $ws.Range("E10") "  -1.65%  "

# Row 11 - TRON
Set-TextValue $ws.Range("D11") "0.07728"
Set-TextValue $ws.Range("E11") "  -0.20%  "

# Row 12 - WrappedEther
Set-TextValue $ws.Range("D12") "1.846.28"
Set-TextValue $ws.Range("E12") "  -2.22%  "

# Row 13 - Polkadot
Set-TextValue $ws.Range("E13") "  -0.71%  "

# Row 14 - Polygon
Set-TextValue $ws.Range("D14") "0.6787"
Set-TextValue $ws.Range("E14") "  -0.46%  "

# Row 15 - ShibaInu
Set-TextValue $ws.Range("D15") "0.00001021"
Set-TextValue $ws.Range("E15") "  -4.91%  "

# Row 16 - Litecoin
Set-TextValue $ws.Range("E16") "  -1.54%  "

# Row 17 - Uniswap
Set-TextValue $ws.Range("D17") "6.144"
Set-TextValue $ws.Range("E17") "  -0.49%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "29.430.02"
Set-TextValue $ws.Range("E18") "  -0.40%  "

# Row 19 - BitcoinCash
Set-TextValue $ws.Range("D19") "228.35"
Set-TextValue $ws.Range("E19") "  -0.12%  "

# Row 20 - Avalanche
Set-TextValue $ws.Range("D20") "12.33"
Set-TextValue $ws.Range("E20") "  -0.10%  "

# Row 21 - Dai
Set-TextValue $ws.Range("D21") "0.9997"
Set-TextValue $ws.Range("E21") "  -0.19%  "

# Row 22 - Chainlink
Set-TextValue $ws.Range("D22") "7.433"
Set-TextValue $ws.Range("E22") "  -0.30%  "

# Row 23 - BinanceUSD
Set-TextValue $ws.Range("D23") "1.0000"
Set-TextValue $ws.Range("E23") "  -0.16%  "

# Row 24 - Monero
Set-TextValue $ws.Range("D24") "159.07"
Set-TextValue $ws.Range("E24") "  +0.66%  "

# Row 25 - Stellar
Set-TextValue $ws.Range("E25") "  -0.40%  "

# Row 26 - Cosmos
Set-TextValue $ws.Range("E26") "  -0.07%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("E27") "  -1.09%  "

# Row 28 - Hedera
Set-TextValue $ws.Range("D28") "0.06280"
Set-TextValue $ws.Range("E28") "  +12.33%  "

# Row 29 - Toncoin
Set-TextValue $ws.Range("D29") "1.386"
Set-TextValue $ws.Range("E29") "  -0.52%  "

# Row 30 - PancakeSwap
Set-TextValue $ws.Range("D30") "1.475"
Set-TextValue $ws.Range("E30") "  +0.87%  "

# Row 31 - Filecoin
Set-TextValue $ws.Range("D31") "4.094"
Set-TextValue $ws.Range("E31") "  -0.99%  "

# Row 32 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D32") "4.059"
Set-TextValue $ws.Range("E32") "  -0.04%  "

# Row 33 - LidoDAOToken
Set-TextValue $ws.Range("E33") "  -1.23%  "

# Row 34 - ARBITRUM
Set-TextValue $ws.Range("D34") "1.142"
Set-TextValue $ws.Range("E34") "  -1.95%  "

# Row 35 - ImmutableX
Set-TextValue $ws.Range("D35") "0.6969"
Set-TextValue $ws.Range("E35") "  -0.08%  "

# Row 36 - HuobiToken
Set-TextValue $ws.Range("D36") "2.582"
Set-TextValue $ws.Range("E36") "  -0.39%  "

# Row 37 - Maker
Set-TextValue $ws.Range("D37") "1.256.55"
Set-TextValue $ws.Range("E37") "  +2.27%  "

# Row 38 - MXToken
Set-TextValue $ws.Range("D38") "2.832"
Set-TextValue $ws.Range("E38") "  +3.75%  "

# Row 39 - VeChain
Set-TextValue $ws.Range("D39") "0.01821"
Set-TextValue $ws.Range("E39") "  +0.89%  "

# Row 40 - FraxShare
Set-TextValue $ws.Range("D40") "6.552"
Set-TextValue $ws.Range("E40") "  +1.79%  "

# Row 41 - TrustWalletToken
Set-TextValue $ws.Range("D41") "0.9095"
Set-TextValue $ws.Range("E41") "  +0.16%  "

# Row 42 - PaxDollar
Set-TextValue $ws.Range("D42") "0.9997"
Set-TextValue $ws.Range("E42") "  -0.19%  "

# Row 43 - RocketPoolETH
Set-TextValue $ws.Range("D43") "2.008.35"
Set-TextValue $ws.Range("E43") "  -14.52%  "

# Row 44 - Quant
Set-TextValue $ws.Range("D44") "101.41"
Set-TextValue $ws.Range("E44") "  -0.81%  "

# Row 45 - Aave
Set-TextValue $ws.Range("D45") "66.34"
Set-TextValue $ws.Range("E45") "  +0.48%  "

# Row 46 & 47 - Aptos/Algorand swap ranking order
Set-TextValue $ws.Range("B46") "Algorand"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D46") "0.1174"
Set-TextValue $ws.Range("E46") "  +2.24%  "

Set-TextValue $ws.Range("B47") "Aptos"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D47") "7.049"
Set-TextValue $ws.Range("E47") "  -2.12%  "

# Row 48 - BabyDogeCoin
Set-TextValue $ws.Range("D48") "0.00000000117"
Set-TextValue $ws.Range("E48") "  -0.09%  "

# Row 49 - EnergySwap
Set-TextValue $ws.Range("D49") "9.048"
Set-TextValue $ws.Range("E49") "  +0.25%  "

# Row 50 - RenderToken
Set-TextValue $ws.Range("D50") "1.685"
Set-TextValue $ws.Range("E50") "  +0.22%  "

# Row 51 - TheSandbox
Set-TextValue $ws.Range("D51") "0.3940"
Set-TextValue $ws.Range("E51") "  -2.08%  "
